$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27

$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "12/21/2025"
$dateCell.ClearFormats()
$ws.Cells.Item($row, 2).Value = 12382.38
$ws.Cells.Item($row, 3).Value = 0.2059261044743488
$ws.Cells.Item($row, 4).Value = 0.7940738955256512
$ws.Cells.Item($row, 5).Value = -128.73
$ws.Cells.Item($row, 6).Value = -26.54
$ws.Cells.Item($row, 7).Value = -20709.23
$ws.Cells.Item($row, 8).Value = -67.81
$ws.Cells.Item($row, 9).Value = -452.99
$ws.Cells.Item($row, 10).Value = -15.09
